$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "GossA-HW45.xpc" to "GossA"
$ws.Name = "GossA"

# Add new row 16 with the 14th HKL entry (reuses the existing
# "HexGrid-60degTilt5degRes" shared string already used in B15)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.004001539351017
$ws.Range("D16").Value = 0.9977992842888757
$ws.Range("E16").Value = 0.9969934576433873
$ws.Range("F16").Value = 0.9965505781855611
$ws.Range("G16").Value = 1.004001539351017
$ws.Range("H16").Value = 0.9977992842888757
$ws.Range("I16").Value = 1.002890434648806
$ws.Range("J16").Value = 1.004183012944848
$ws.Range("K16").Value = 1.001764705882353
$ws.Range("L16").Value = 0.9988381534760072
$ws.Range("M16").Value = 1.004001539351017
$ws.Range("N16").Value = 0.9973963709661315
$ws.Range("O16").Value = 0.9988362148672103
$ws.Range("P16").Value = 1.000377645802607

# Match style of A16 with A15 (bordered/centered header-like style)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
